# Apply "Updated RPS with current policy values" changes.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "About" sheet: refresh the Sources section with a real citation.
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# The two old BAU/no-target note rows (old rows 10 & 11) go away entirely.
# Delete from the bottom up so row numbers above stay stable.
$wsAbout.Rows.Item(11).Delete() | Out-Null
$wsAbout.Rows.Item(10).Delete() | Out-Null

# Rows 7-9 had a redundant "apply default font" style on them; drop it.
$wsAbout.Range("A7:A9").ClearFormats() | Out-Null

# Add the NREL source link first (so its shared-string slot precedes the
# note text below, matching insertion order).
$wsAbout.Hyperlinks.Add($wsAbout.Range("B5"), "https://www.nrel.gov/docs/fy22osti/82580.pdf") | Out-Null

# Replace the old "no RPS target" note with the new sourced note.
$rsquo = [char]0x2019
$newNote = "Mexico" + $rsquo + "s energy transition law established a target for meeting at" + "`n" + "least 35% of its electricity generation from clean energy sources" + "`n" + "by 2024. "
$wsAbout.Range("B3").Value2 = $newNote
$wsAbout.Range("B3").WrapText = $true
$wsAbout.Rows.Item(3).RowHeight = 48

$wsAbout.Range("B10").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) "BRPSPTY" sheet: the BAU RPS values move from 0 to the current 35%
#    clean-energy-by-2024 policy figure.
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("BRPSPTY")

# D1 ("Year2017" header) had the same redundant default-font style; clear it.
$wsData.Range("D1").ClearFormats() | Out-Null

$wsData.Rows.Item(2).RowHeight = 32
$wsData.Range("B2:CG2").Value = 0.35

$wsData.Activate() | Out-Null
$wsData.Range("B2:CG2").Select() | Out-Null
